$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: 中国南方电网 / http://www.bidding.csg.cn/zbgg/index.jhtml ---
$ws.Range("A16").Value = "中国南方电网"
$ws.Range("B16").Value = "http://www.bidding.csg.cn/zbgg/index.jhtml"

# --- Row 17: 中国石化 / https://ebidding.sinopec.com/TPWeb4AAA/jyxx/002002/ ---
$ws.Range("A17").Value = "中国石化"
$ws.Range("B17").Value = "https://ebidding.sinopec.com/TPWeb4AAA/jyxx/002002/"

# --- Row 18: 中国石油 (no link/value in column B) ---
$ws.Range("A18").Value = "中国石油"

# Match formatting of the existing "category" rows (A10/C10 use the
# bordered, centered look shared by A2:A4, A6:A7, A8:A9, A11:A15, A5/A10, ...)
$ws.Range("A10:C10").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)
$ws.Range("A17:C17").PasteSpecial(-4122)
$ws.Range("A18:C18").PasteSpecial(-4122)

# Match formatting of the existing hyperlink-text cells in column B
$ws.Range("B2").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B18").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row heights to match the rest of the table
$ws.Rows.Item(16).RowHeight = 16.5
$ws.Rows.Item(17).RowHeight = 16.5
$ws.Rows.Item(18).RowHeight = 16.5

# Hyperlinks for the two new URL cells
$ws.Hyperlinks.Add($ws.Range("B16"), "http://www.bidding.csg.cn/zbgg/index.jhtml")
$ws.Hyperlinks.Add($ws.Range("B17"), "https://ebidding.sinopec.com/TPWeb4AAA/jyxx/002002/")

# Re-apply the look for the hyperlink cells after Hyperlinks.Add resets it
$ws.Range("B2").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B17").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# B18 stays blank, but still carries the same look as B16/B17
$ws.Range("B2").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = ""
$ws.Application.CutCopyMode = $false

$ws.Range("B18").Select()
